# "change name file part"
# Renames the uploaded test part files referenced on each order sheet:
#   "stp-testing 2.stp"   -> "stp-testing-2.stp"
#   "step-testing 2.step" -> "step-testing-2.step"
# (the space before the "2" is replaced with a hyphen)

$wb = $excel.ActiveWorkbook

# Sheets FPA001, FPA002-003-005-007, FPA004-006-010, FPA008-009, BTMI003, BTMI015
# each keep the "stp-testing 2.stp" file name in cell A3.
$stpSheets = @("FPA002-003-005-007", "FPA004-006-010", "FPA008-009", "BTMI003", "BTMI015")
foreach ($sheetName in $stpSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A3").Value = "stp-testing-2.stp"
}

# Sheet BTMI002 keeps the "step-testing 2.step" file name in cell A4.
$ws5 = $wb.Worksheets.Item("BTMI002")
$ws5.Activate()
$ws5.Range("A4").Value = "step-testing-2.step"
$ws5.Range("A4").Select() | Out-Null

# Sheet FPA001 (the first/active sheet) keeps the "stp-testing 2.stp" file
# name in cell A3 as well - edit it last so it ends up the active tab again.
$ws1 = $wb.Worksheets.Item("FPA001")
$ws1.Activate()
$ws1.Range("A3").Value = "stp-testing-2.stp"
